$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New regenerated s_val data (filtered save games), columns B:G, rows 2-8
$data = @(
    @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 1, 15.88780690183548),
    @(0.2917716402565462, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 0, 1.579467928156517),
    @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 10.19245300693656, 0, 13.45301510845117),
    @(0.6606524410359556, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 1, 2.214453472130288),
    @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548),
    @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 3.754798637575387),
    @(3.286832544864788, 1.655778082260271, 3.537761648806719, 10.19245300693656, 0, 18.67282528286833)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
